$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp on row 13 (tiny precision fix)
$ws.Range("A13").Value = 45878.54183790509

# Append new row 14 with the latest sensor reading
$ws.Range("A14").Value = 45878.58356094634
$ws.Range("B14").Value = 2025
$ws.Range("C14").Value = 37
$ws.Range("D14").Value = 18.94
$ws.Range("E14").Value = 77.43000000000001
$ws.Range("F14").Value = 94.92
$ws.Range("G14").Value = 13.84
$ws.Range("H14").Value = "ESE"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "14:00:19"

# Match the date/time number format + style used by the rest of column A
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
